$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at CK (89th column). This pushes the existing
# "nom" (CK) and "url_produit" (CL) columns one to the right, to CL and CM.
$ws.Columns("CK").Insert()

# The newly inserted CK column becomes the latest price-history snapshot,
# a duplicate of the previous last snapshot column (now at CJ) for every
# data row.
$src = $ws.Range("CJ2:CJ206")
$dst = $ws.Range("CK2:CK206")
$dst.Value = $src.Value()

# Header gets the new snapshot's timestamp.
$ws.Range("CK1").Value = "2026-01-31 17:15:20"
